$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $r1, $r2, $firstCol, $lastCol) {
    $col = $firstCol
    while ($col -le $lastCol) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)
        $v1 = $cell1.Value()
        $v2 = $cell2.Value()
        $cell1.Value = $v2
        $cell2.Value = $v1
        $col = $col + 1
    }
}

# Column B = 2 ... Column AC = 29
$firstCol = 2
$lastCol = 29

# Swap rows 108 and 110 (B..AC)
Swap-Rows $ws 108 110 $firstCol $lastCol

# Swap rows 112 and 114 (B..AC)
Swap-Rows $ws 112 114 $firstCol $lastCol

# Swap rows 137 and 138 (B..AC)
Swap-Rows $ws 137 138 $firstCol $lastCol

# Swap rows 139 and 140 (B..AC)
Swap-Rows $ws 139 140 $firstCol $lastCol

# Direct value updates on row 173
$ws.Cells.Item(173, 14).Value = 2.375   # N173 oddH
$ws.Cells.Item(173, 16).Value = 3.3     # P173 oddA
$ws.Cells.Item(173, 18).Value = 1.95    # R173 oddAHH
$ws.Cells.Item(173, 19).Value = 1.9     # S173 oddAHA
$ws.Cells.Item(173, 21).Value = 1.95    # U173 oddAHOver
$ws.Cells.Item(173, 22).Value = 1.9     # V173 oddAHUnder
